$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.757.24'
$ws.Range('E2').Value = '  -5.89%  '
$ws.Range('D3').Value = '2.978.51'
$ws.Range('E3').Value = '  -6.24%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.50'
$ws.Range('E5').Value = '  -4.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '124.55'
$ws.Range('E6').Value = '  -8.63%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = '2.972.48'
$ws.Range('E8').Value = '  -6.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.498'
$ws.Range('E9').Value = '  -2.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.129'
$ws.Range('E10').Value = '  -9.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.07'
$ws.Range('E11').Value = '  -5.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.436'
$ws.Range('E12').Value = '  -4.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000217'
$ws.Range('E13').Value = '  -9.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.34'
$ws.Range('E14').Value = '  -6.90%  '
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').Value = '3.466.28'
$ws.Range('E16').Value = '  -6.22%  '
$ws.Range('D17').Value = '2.986.58'
$ws.Range('E17').Value = '  -5.95%  '
$ws.Range('D18').Value = '59.875.34'
$ws.Range('E18').Value = '  -5.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.37'
$ws.Range('E19').Value = '  -2.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '424.43'
$ws.Range('E20').Value = '  -7.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.02'
$ws.Range('E21').Value = '  -6.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.662'
$ws.Range('E22').Value = '  -4.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.00'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.95'
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '78.92'
$ws.Range('E25').Value = '  -5.00%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.52'
$ws.Range('E28').Value = '  -6.02%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.93'
$ws.Range('E29').Value = '  -7.52%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.21'
$ws.Range('E30').Value = '  -7.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.11'
$ws.Range('E31').Value = '  -10.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.03'
$ws.Range('E32').Value = '  -9.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0930'
$ws.Range('E33').Value = '  -8.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.935'
$ws.Range('E34').Value = '  -8.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.57'
$ws.Range('E35').Value = '  -5.20%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '50.12'
$ws.Range('E36').Value = '  -2.48%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.03'
$ws.Range('E37').Value = '  -17.70%  '
$ws.Range('D38').Value = '0.0₃0654'
$ws.Range('E38').Value = '  -10.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.29'
$ws.Range('E39').Value = '  +2.16%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0352'
$ws.Range('E40').Value = '  -10.26%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.108'
$ws.Range('E41').Value = '  -5.02%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '378.45'
$ws.Range('E42').Value = '  -4.78%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.653.47'
$ws.Range('E43').Value = '  -5.48%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.45'
$ws.Range('E44').Value = '  -8.44%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '119.91'
$ws.Range('E46').Value = '  -7.45%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.233'
$ws.Range('E47').Value = '  -8.14%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.98'
$ws.Range('E48').Value = '  -7.68%  '
$ws.Range('E49').Value = '  -4.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.41'
$ws.Range('E50').Value = '  -8.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.131'
$ws.Range('E51').Value = '  -1.21%  '
